$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1.010622732913008
$ws.Range("D2").Value = 1.013268254309555
$ws.Range("E2").Value = 1.012829630966766
$ws.Range("F2").Value = 1.008866401775824
$ws.Range("J2").Value = 1.015875912179003
$ws.Range("K2").Value = 1.016129916656902
$ws.Range("L2").Value = 1.015692611750708
$ws.Range("M2").Value = 1.011741352102767
$ws.Range("N2").Value = 1.017318572496894
$ws.Range("C3").Value = 1.012267720130226
$ws.Range("D3").Value = 1.014824216018873
$ws.Range("E3").Value = 1.014248996179511
$ws.Range("F3").Value = 1.011170520799082
$ws.Range("J3").Value = 1.017150088860458
$ws.Range("K3").Value = 1.017488479125081
$ws.Range("L3").Value = 1.016914860441044
$ws.Range("M3").Value = 1.013844987741668
$ws.Range("N3").Value = 1.018594558655386
$ws.Range("C4").Value = 1.013329534166127
$ws.Range("D4").Value = 1.01582882617506
$ws.Range("E4").Value = 1.015165307789248
$ws.Range("F4").Value = 1.012658018727509
$ws.Range("J4").Value = 1.017971744485487
$ws.Range("K4").Value = 1.018364895301193
$ws.Range("L4").Value = 1.017703128852679
$ws.Range("M4").Value = 1.015202482377503
$ws.Range("N4").Value = 1.019417381125647
$ws.Range("C5").Value = 1.013775313124032
$ws.Range("D5").Value = 1.016250651085926
$ws.Range("E5").Value = 1.015550031501987
$ws.Range("F5").Value = 1.013282571957313
$ws.Range("J5").Value = 1.018316505173811
$ws.Range("K5").Value = 1.018732715462581
$ws.Range("L5").Value = 1.018033903713063
$ws.Range("M5").Value = 1.015772313472751
$ws.Range("N5").Value = 1.019762631413693
$ws.Range("C6").Value = 1.013850126184366
$ws.Range("D6").Value = 1.016321447647368
$ws.Range("E6").Value = 1.015614599704553
$ws.Range("F6").Value = 1.0133873916259
$ws.Range("J6").Value = 1.018374353465014
$ws.Range("K6").Value = 1.018794437791005
$ws.Range("L6").Value = 1.01808940663465
$ws.Range("M6").Value = 1.015867940953244
$ws.Range("N6").Value = 1.019820561856106
$ws.Range("C7").Value = 1.013335493053048
$ws.Range("D7").Value = 1.015834464617905
$ws.Range("E7").Value = 1.015170450409669
$ws.Range("F7").Value = 1.012666367108049
$ws.Range("J7").Value = 1.017976353784442
$ws.Range("K7").Value = 1.018369812571446
$ws.Range("L7").Value = 1.017707551076573
$ws.Range("M7").Value = 1.015210099831562
$ws.Range("N7").Value = 1.019421996970336
$ws.Range("C8").Value = 1.011179212277533
$ws.Range("D8").Value = 1.013794562926315
$ws.Range("E8").Value = 1.013309756813115
$ws.Range("F8").Value = 1.009645814975661
$ws.Range("J8").Value = 1.016307117841988
$ws.Range("K8").Value = 1.016589608609779
$ws.Range("L8").Value = 1.016106223262074
$ws.Range("M8").Value = 1.012453066933475
$ws.Range("N8").Value = 1.017750390521379
$ws.Range("C9").Value = 1.007358940637754
$ws.Range("D9").Value = 1.010182517007349
$ws.Range("E9").Value = 1.01001425211909
$ws.Range("F9").Value = 1.004295744692176
$ws.Range("J9").Value = 1.013343546487522
$ws.Range("K9").Value = 1.013431680646142
$ws.Range("L9").Value = 1.013263996046943
$ws.Range("M9").Value = 1.007565337116161
$ws.Range("N9").Value = 1.014782610555663
$ws.Range("C10").Value = 1.004797244699635
$ws.Range("D10").Value = 1.007761862039963
$ws.Range("E10").Value = 1.007805237397172
$ws.Range("F10").Value = 1.000708777428092
$ws.Range("J10").Value = 1.011352153591791
$ws.Range("K10").Value = 1.011311481969041
$ws.Range("L10").Value = 1.011354692591863
$ws.Range("M10").Value = 1.00428539486348
$ws.Range("N10").Value = 1.012788389653605
$ws.Range("C11").Value = 1.003684261183982
$ws.Range("D11").Value = 1.006710501723597
$ws.Range("E11").Value = 1.006845688626189
$ws.Range("F11").Value = 0.9991503778886717
$ws.Range("J11").Value = 1.010485965543977
$ws.Range("K11").Value = 1.010389694353908
$ws.Range("L11").Value = 1.0105243472595
$ws.Range("M11").Value = 1.002859689474237
$ws.Range("N11").Value = 1.01192097151941
$ws.Range("C12").Value = 1.00327026682997
$ws.Range("D12").Value = 1.006319481849011
$ws.Range("E12").Value = 1.006488799212478
$ws.Range("F12").Value = 0.9985707002937256
$ws.Range("J12").Value = 1.010163623420673
$ws.Range("K12").Value = 1.010046725662995
$ws.Range("L12").Value = 1.010215364615088
$ws.Range("M12").Value = 1.002329265808783
$ws.Range("N12").Value = 1.011598171633318
$ws.Range("C13").Value = 1.003359096700946
$ws.Range("D13").Value = 1.00640337974255
$ws.Range("E13").Value = 1.006565374725978
$ws.Range("F13").Value = 0.9986950807090053
$ws.Range("J13").Value = 1.010232794355919
$ws.Range("K13").Value = 1.010120319918557
$ws.Range("L13").Value = 1.010281667800164
$ws.Range("M13").Value = 1.00244308260659
$ws.Range("N13").Value = 1.011667440799226
$ws.Range("C14").Value = 1.003650052279565
$ws.Range("D14").Value = 1.006678190153546
$ws.Range("E14").Value = 1.006816197716495
$ws.Range("F14").Value = 0.9991024784588273
$ws.Range("J14").Value = 1.010459332998334
$ws.Range("K14").Value = 1.010361356291173
$ws.Range("L14").Value = 1.010498818069175
$ws.Range("M14").Value = 1.002815862078011
$ws.Range("N14").Value = 1.011894301152497
$ws.Range("C15").Value = 1.003829241826
$ws.Range("D15").Value = 1.006847443437462
$ws.Range("E15").Value = 1.006970675272938
$ws.Range("F15").Value = 0.999353380142539
$ws.Range("J15").Value = 1.010598830796842
$ws.Range("K15").Value = 1.010509789992769
$ws.Range("L15").Value = 1.01063253749945
$ws.Range("M15").Value = 1.003045429803249
$ws.Range("N15").Value = 1.01203399705388
$ws.Range("C16").Value = 1.004871029079135
$ws.Range("D16").Value = 1.007831568482542
$ws.Range("E16").Value = 1.007868854363161
$ws.Range("F16").Value = 1.000812090378618
$ws.Range("J16").Value = 1.01140955614799
$ws.Range("K16").Value = 1.01137257810577
$ws.Range("L16").Value = 1.011409722754908
$ws.Range("M16").Value = 1.004379896328309
$ws.Range("N16").Value = 1.012845873728017
$ws.Range("C17").Value = 1.00552349678381
$ws.Range("D17").Value = 1.008448015285461
$ws.Range("E17").Value = 1.008431436462027
$ws.Range("F17").Value = 1.001725679268547
$ws.Range("J17").Value = 1.011917047514307
$ws.Range("K17").Value = 1.01191277329278
$ws.Range("L17").Value = 1.011896255789869
$ws.Range("M17").Value = 1.005215485621803
$ws.Range("N17").Value = 1.013354085790285
$ws.Range("C18").Value = 1.00590370941189
$ws.Range("D18").Value = 1.00880727060816
$ws.Range("E18").Value = 1.008759289670581
$ws.Range("F18").Value = 1.00225805900373
$ws.Range("J18").Value = 1.012212683207891
$ws.Range("K18").Value = 1.012227501417766
$ws.Range("L18").Value = 1.012179695460119
$ws.Range("M18").Value = 1.005702345219156
$ws.Range("N18").Value = 1.013650141320464
$ws.Range("C19").Value = 1.006033291474621
$ws.Range("D19").Value = 1.008929715713357
$ws.Range("E19").Value = 1.008871030162452
$ws.Range("F19").Value = 1.00243950275952
$ws.Range("J19").Value = 1.012313424113663
$ws.Range("K19").Value = 1.012334755224156
$ws.Range("L19").Value = 1.012276282654582
$ws.Range("M19").Value = 1.005868263575752
$ws.Range("N19").Value = 1.013751025289877
$ws.Range("C20").Value = 1.005453530624437
$ws.Range("D20").Value = 1.008381908323378
$ws.Range("E20").Value = 1.008371106976772
$ws.Range("F20").Value = 1.001627711917561
$ws.Range("J20").Value = 1.011862637435466
$ws.Range("K20").Value = 1.011854852665216
$ws.Range("L20").Value = 1.011844091372952
$ws.Range("M20").Value = 1.00512588938099
$ws.Range("N20").Value = 1.013299598442891
$ws.Range("C21").Value = 1.00356438935315
$ws.Range("D21").Value = 1.006597279217181
$ws.Range("E21").Value = 1.006742349700188
$ws.Range("F21").Value = 0.9989825328530684
$ws.Range("J21").Value = 1.010392639736102
$ws.Range("K21").Value = 1.01029039310712
$ws.Range("L21").Value = 1.010434888172092
$ws.Range("M21").Value = 1.002706111623407
$ws.Range("N21").Value = 1.011827513178183
$ws.Range("C22").Value = 1.002373230232304
$ws.Range("D22").Value = 1.005472323144594
$ws.Range("E22").Value = 1.00571555694255
$ws.Range("F22").Value = 0.9973146500667869
$ws.Range("J22").Value = 1.00946490717971
$ws.Range("K22").Value = 1.009303416377553
$ws.Range("L22").Value = 1.009545645695087
$ws.Range("M22").Value = 1.00117974896345
$ws.Range("N22").Value = 1.010898463135149
$ws.Range("C23").Value = 1.003005013103309
$ws.Range("D23").Value = 1.0060689632019
$ws.Range("E23").Value = 1.006260142868672
$ws.Range("F23").Value = 0.9981992885366436
$ws.Range("J23").Value = 1.009957051354936
$ws.Range("K23").Value = 1.009826953201411
$ws.Range("L23").Value = 1.010017359913397
$ws.Range("M23").Value = 1.001989382774208
$ws.Range("N23").Value = 1.011391306211553
$ws.Range("C24").Value = 1.005485146434493
$ws.Range("D24").Value = 1.008411780164376
$ws.Range("E24").Value = 1.008398368174521
$ws.Range("F24").Value = 1.001671980694063
$ws.Range("J24").Value = 1.011887224167814
$ws.Range("K24").Value = 1.011881025611366
$ws.Range("L24").Value = 1.011867663298968
$ws.Range("M24").Value = 1.005166375688318
$ws.Range("N24").Value = 1.013324220091219
$ws.Range("C25").Value = 1.008349117105379
$ws.Range("D25").Value = 1.011118477954688
$ws.Range("E25").Value = 1.010868279412309
$ws.Range("F25").Value = 1.005682303827192
$ws.Range("J25").Value = 1.01411240667111
$ws.Range("K25").Value = 1.014250649547849
$ws.Range("L25").Value = 1.01400128018193
$ws.Range("M25").Value = 1.00883259329417
$ws.Range("N25").Value = 1.015552562608901
